$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q3" sheet right before "2022-Q2", by duplicating the
#    "2022-Q2" sheet (so it inherits identical formatting/styles) and then
#    overwriting its contents with the Q3 fund-holding data.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Make sure column A has the right style as far down as row 10 (2022-Q2 only
# had data through row 7, but 2022-Q3 needs rows 2-10).
$q3.Cells.Item(2,1).Copy($q3.Cells.Item(8,1))
$q3.Cells.Item(2,1).Copy($q3.Cells.Item(9,1))
$q3.Cells.Item(2,1).Copy($q3.Cells.Item(10,1))

$q3Rows = @(
    @(0, "450008", "国富沪深300指数增强",           "4.26", "91.78", "1.99", "0.0848", 9),
    @(1, "002849", "金信智能中国2025灵活配置混合",   "1.09", "78.83", "5.04", "0.0549", 5),
    @(2, "510760", "国泰上证综合ETF",                "2.75", "94.96", "1.37", "0.0377", 6),
    @(3, "004730", "建信量化事件驱动股票",           "0.46", "81.67", "1.67", "0.0077", 9),
    @(4, "004892", "华润元大成长精选股票C",          "0.12", "94.00", "1.89", "0.0023", 8),
    @(5, "005053", "银河量化价值混合A",              "0.10", "78.55", "1.74", "0.0017", 5),
    @(6, "005126", "银河量化稳进混合",               "0.13", "55.69", "1.09", "0.0014", 7),
    @(7, "004891", "华润元大成长精选股票A",          "0.04", "94.00", "1.89", "0.0008", 8),
    @(8, "013026", "银河量化价值混合C",              "0.00", "78.55", "1.74", 0,      5)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Cells.Item($r,1).Value = $row[0]

    # Columns B/D/E/F/G hold numeric-looking text (fund codes, percentages,
    # etc.) that must stay text, not get auto-coerced to Number. A leading
    # apostrophe forces text entry; resetting the style back to "Normal"
    # afterwards strips the quote-prefix formatting Excel applies so the
    # cell ends up with plain default styling (matching the source data,
    # which was written as plain inline strings, not typed via the UI).
    $q3.Cells.Item($r,2).Value = "'" + $row[1]
    $q3.Cells.Item($r,2).Style = "Normal"

    $q3.Cells.Item($r,3).Value = $row[2]

    $q3.Cells.Item($r,4).Value = "'" + $row[3]
    $q3.Cells.Item($r,4).Style = "Normal"

    $q3.Cells.Item($r,5).Value = "'" + $row[4]
    $q3.Cells.Item($r,5).Style = "Normal"

    $q3.Cells.Item($r,6).Value = "'" + $row[5]
    $q3.Cells.Item($r,6).Style = "Normal"

    if ($r -eq 10) {
        $q3.Cells.Item($r,7).Value = $row[6]
    } else {
        $q3.Cells.Item($r,7).Value = "'" + $row[6]
        $q3.Cells.Item($r,7).Style = "Normal"
    }

    $q3.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: shift every existing quarter row down
#    by one, insert the new 2022-Q3 totals at row 2, and append the row that
#    drops off the bottom (2020-Q4) as the new last row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend column-A styling down to the new row 9.
$total.Cells.Item(2,1).Copy($total.Cells.Item(9,1))

$totalRows = @(
    @("2022-Q3", 9,  0.19),
    @("2022-Q2", 6,  0.14),
    @("2022-Q1", 15, 1.81),
    @("2021-Q4", 6,  0.32),
    @("2021-Q3", 2,  0.13),
    @("2021-Q2", 4,  0.3),
    @("2021-Q1", 6,  0.24),
    @("2020-Q4", 6,  0.09)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r,1).Value = $r - 2
    $total.Cells.Item($r,2).Value = $row[0]
    $total.Cells.Item($r,3).Value = $row[1]
    $total.Cells.Item($r,4).Value = $row[2]
    $r = $r + 1
}
